$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The log sheet gets one new row appended at the bottom (row 52), following
# the same pattern/style as the previous row (row 51). Copy formatting first
# so the new row keeps the same cell style (centered alignment, etc.), then
# fill in the values for the new run-log entry.
$ws.Range("A51:H51").Copy()
$ws.Range("A52:H52").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(52, 1).Value = "2025-08-24 09:36:25 UTC"
$ws.Cells.Item(52, 2).Value = "2025-08-24 15:06:25 IST"
$ws.Cells.Item(52, 3).Value = "SKIPPED"
$ws.Cells.Item(52, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item(52, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item(52, 6).Value = ""
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = ""
